$wb = $excel.ActiveWorkbook

# "Overview" sheet: the Latest HO Xliff Generate Date for the 348fa1fa file
# is refreshed to the new handback-generation timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 16:57:46"

# "zh-cn" sheet: refresh handoff/handback datetimes for the 348fa1fa row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 16:57:42"
$wsZhCn.Range("K2").Value = "2016-09-01 16:57:59"

# "de-de" sheet: refresh handoff/handback datetimes for the 348fa1fa row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 16:57:46"
$wsDeDe.Range("K2").Value = "2016-09-01 16:58:18"
